# Refresh the hourly crypto price/volume snapshot (GitHub Actions bot run).
#
# Column D ("Price") cells are stored as plain text in the workbook (values
# like "35.950.81" or "0.0800" aren't valid numbers / would lose their
# trailing zeros and literal dot-grouping if Excel auto-coerced them to
# numeric). To keep them as exact text - same as the original - each D cell
# is switched to the Text number format before the write, then has its
# format cleared again afterwards so no stray style is left referenced on
# the cell (matches the unstyled inlineStr cells in the source file).
#
# Column E ("Volume(1h)") cells are padded percent strings (e.g.
# "  -4.59%  ") which Excel already keeps as text, so a plain .Value
# assignment is enough there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '35.950.81'
$c.ClearFormats()
$ws.Range("E2").Value = '  -4.59%  '
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.956.59'
$c.ClearFormats()
$ws.Range("E3").Value = '  -4.21%  '
$ws.Range("E4").Value = '  +0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '240.86'
$c.ClearFormats()
$ws.Range("E5").Value = '  -4.06%  '
$ws.Range("E6").Value = '  -3.94%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '60.22'
$c.ClearFormats()
$ws.Range("E7").Value = '  -8.50%  '
$ws.Range("E8").Value = '  +0.05%  '
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.373'
$c.ClearFormats()
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E10").Value = '  -5.39%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.0800'
$c.ClearFormats()
$ws.Range("E11").Value = '  +6.23%  '
$ws.Range("E12").Value = '  -1.47%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.850'
$c.ClearFormats()
$ws.Range("E13").Value = '  -5.54%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '22.00'
$c.ClearFormats()
$ws.Range("E14").Value = '  +7.19%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '13.94'
$c.ClearFormats()
$ws.Range("E15").Value = '  -7.31%  '
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '2.239.95'
$c.ClearFormats()
$ws.Range("E16").Value = '  -4.32%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '5.38'
$c.ClearFormats()
$ws.Range("E17").Value = '  -3.32%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '1.962.32'
$c.ClearFormats()
$ws.Range("E18").Value = '  -3.96%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '35.874.75'
$c.ClearFormats()
$ws.Range("E19").Value = '  -4.46%  '
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '70.69'
$c.ClearFormats()
$ws.Range("E20").Value = '  -3.44%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '0.0₃0854'
$c.ClearFormats()
$ws.Range("E21").Value = '  -2.07%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '234.76'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.96%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '5.18'
$c.ClearFormats()
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("E24").Value = '  +0.19%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.52'
$c.ClearFormats()
$ws.Range("E25").Value = '  -5.80%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.ClearFormats()
$ws.Range("E26").Value = '  -4.26%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '9.69'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.50%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '159.83'
$c.ClearFormats()
$ws.Range("E28").Value = '  -2.87%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '19.64'
$c.ClearFormats()
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("E30").Value = '  +14.54%  '
$ws.Range("E31").Value = '  -2.07%  '
$ws.Range("E32").Value = '  -6.99%  '
$ws.Range("E33").Value = '  -6.55%  '
$ws.Range("E34").Value = '  +0.62%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '4.38'
$c.ClearFormats()
$ws.Range("E35").Value = '  -7.01%  '
$ws.Range("E36").Value = '  +3.08%  '
$ws.Range("E37").Value = '  +0.18%  '
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '2.25'
$c.ClearFormats()
$ws.Range("E38").Value = '  -7.77%  '
$ws.Range("E39").Value = '  -0.88%  '
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '3.04'
$c.ClearFormats()
$ws.Range("E40").Value = '  +10.12%  '
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.0984'
$c.ClearFormats()
$ws.Range("E41").Value = '  -4.35%  '
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '1.22'
$c.ClearFormats()
$ws.Range("E42").Value = '  -0.66%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '2.85'
$c.ClearFormats()
$ws.Range("E43").Value = '  -3.14%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.0211'
$c.ClearFormats()
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("E45").Value = '  -4.89%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '91.38'
$c.ClearFormats()
$ws.Range("E46").Value = '  -3.59%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '15.91'
$c.ClearFormats()
$ws.Range("E47").Value = '  -5.36%  '
$ws.Range("E48").Value = '  -7.64%  '
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '1.329.87'
$c.ClearFormats()
$ws.Range("E49").Value = '  -6.51%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '2.82'
$c.ClearFormats()
$ws.Range("E50").Value = '  -4.08%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '2.137.33'
$c.ClearFormats()
$ws.Range("E51").Value = '  -4.08%  '
